$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 750.0909
$ws.Range("I4").Value = 527.6667
$ws.Range("J4").Value = 1751
$ws.Range("K4").Value = 527.6667
$ws.Range("L4").Value = 1751
$ws.Range("M4").Value = -413.6667
$ws.Range("N4").Value = -1979
$ws.Range("H18").Value = 729
$ws.Range("I18").Value = 729
$ws.Range("K18").Value = 729
$ws.Range("M18").Value = -445
$ws.Range("H53").Value = 1223.1765
$ws.Range("J53").Value = 859.2727
$ws.Range("L53").Value = 859.2727
$ws.Range("N53").Value = -2133.2727
$ws.Range("H125").Value = 5879.25
$ws.Range("I125").Value = 4999.6665
$ws.Range("J125").Value = 6407
$ws.Range("K125").Value = 44996.9985
$ws.Range("L125").Value = 57663
$ws.Range("M125").Value = -42536.9985
$ws.Range("N125").Value = -62583
$ws.Range("H137").Value = 1718.375
$ws.Range("I137").Value = 1596.9445
$ws.Range("J137").Value = 2082.6667
$ws.Range("K137").Value = 4790.833500000001
$ws.Range("L137").Value = 6248.000100000001
$ws.Range("M137").Value = -2240.833500000001
$ws.Range("N137").Value = -11348.0001
$ws.Range("H138").Value = 4579.9263
$ws.Range("J138").Value = 3902.4338
$ws.Range("L138").Value = 11707.3014
$ws.Range("N138").Value = -21987.3014

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1146.381
$ws.Range("I2").Value = 954.05884
$ws.Range("K2").Value = 954.05884
$ws.Range("M2").Value = -841.05884
$ws.Range("H28").Value = 30564.857
$ws.Range("I28").Value = 30564.857
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 30564.857
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -30372.857
$ws.Range("H32").Value = 9281.234
$ws.Range("J32").Value = 30310
$ws.Range("L32").Value = 30310
$ws.Range("N32").Value = -30884
$ws.Range("H45").Value = 2630.0908
$ws.Range("I45").Value = 1285.3334
$ws.Range("K45").Value = 1285.3334
$ws.Range("M45").Value = -908.3334
$ws.Range("H74").Value = 1896.375
$ws.Range("I74").Value = 1429.3793
$ws.Range("K74").Value = 1429.3793
$ws.Range("M74").Value = -555.3793000000001
$ws.Range("H77").Value = 1896.375
$ws.Range("I77").Value = 1429.3793
$ws.Range("K77").Value = 7146.896500000001
$ws.Range("M77").Value = -2778.896500000001
$ws.Range("H99").Value = 30564.857
$ws.Range("I99").Value = 30564.857
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 30564.857
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -27569.857
$ws.Range("H108").Value = 100000
$ws.Range("J108").Value = 100000
$ws.Range("L108").Value = 100000
$ws.Range("N108").Value = -107680
$ws.Range("H116").Value = 1146.381
$ws.Range("I116").Value = 954.05884
$ws.Range("K116").Value = 954.05884
$ws.Range("M116").Value = 1339.94116
$ws.Range("H122").Value = 8225.963
$ws.Range("I122").Value = 8084.16
$ws.Range("K122").Value = 24252.48
$ws.Range("M122").Value = -21802.48

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1146.381
$ws.Range("I3").Value = 954.05884
$ws.Range("K3").Value = 954.05884
$ws.Range("M3").Value = -840.05884
$ws.Range("H64").Value = 1361.6666
$ws.Range("I64").Value = 1088
$ws.Range("K64").Value = 1088
$ws.Range("M64").Value = -863
$ws.Range("H67").Value = 1361.6666
$ws.Range("I67").Value = 1088
$ws.Range("K67").Value = 1088
$ws.Range("M67").Value = -308

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 418
$ws.Range("I7").Value = 435
$ws.Range("J7").Value = 401
$ws.Range("K7").Value = 435
$ws.Range("L7").Value = 401
$ws.Range("M7").Value = -322
$ws.Range("N7").Value = -627
$ws.Range("H31").Value = 48386.816
$ws.Range("I31").Value = 3251.0833
$ws.Range("J31").Value = 102549.7
$ws.Range("K31").Value = 3251.0833
$ws.Range("L31").Value = 102549.7
$ws.Range("M31").Value = -2956.0833
$ws.Range("N31").Value = -103139.7
$ws.Range("H34").Value = 48386.816
$ws.Range("I34").Value = 3251.0833
$ws.Range("J34").Value = 102549.7
$ws.Range("K34").Value = 3251.0833
$ws.Range("L34").Value = 102549.7
$ws.Range("M34").Value = -3049.0833
$ws.Range("N34").Value = -102953.7
$ws.Range("H141").Value = 422359.72
$ws.Range("J141").Value = 676481.5
$ws.Range("L141").Value = 676481.5
$ws.Range("N141").Value = -686841.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1553.9
$ws.Range("J92").Value = 1593.3334
$ws.Range("L92").Value = 4780.0002
$ws.Range("N92").Value = -7276.0002
$ws.Range("H97").Value = 599.5
$ws.Range("J97").Value = 599.5
$ws.Range("L97").Value = 1798.5
$ws.Range("N97").Value = -2790.5
$ws.Range("H121").Value = 54970.61
$ws.Range("I121").Value = 299.6
$ws.Range("J121").Value = 70157
$ws.Range("K121").Value = 898.8000000000001
$ws.Range("L121").Value = 210471
$ws.Range("M121").Value = 411.1999999999999
$ws.Range("N121").Value = -213091
$ws.Range("H131").Value = 7838092
$ws.Range("I131").Value = 14014676
$ws.Range("K131").Value = 42044028
$ws.Range("M131").Value = -42038988

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 348.57693
$ws.Range("J2").Value = 887.2222
$ws.Range("L2").Value = 887.2222
$ws.Range("N2").Value = -1113.2222
$ws.Range("H102").Value = 17569.354
$ws.Range("I102").Value = 19179.334
$ws.Range("K102").Value = 19179.334
$ws.Range("M102").Value = -17557.334
$ws.Range("H122").Value = 43320.08
$ws.Range("I122").Value = 66666.5
$ws.Range("K122").Value = 199999.5
$ws.Range("M122").Value = -197549.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3245.5334
$ws.Range("I40").Value = 2719.1
$ws.Range("J40").Value = 4298.4
$ws.Range("K40").Value = 2719.1
$ws.Range("L40").Value = 4298.4
$ws.Range("M40").Value = -2583.1
$ws.Range("N40").Value = -4570.4
$ws.Range("H55").Value = 1325.8667
$ws.Range("I55").Value = 322.25
$ws.Range("K55").Value = 322.25
$ws.Range("M55").Value = -149.25
$ws.Range("H61").Value = 2727
$ws.Range("I61").Value = 2727
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2727
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2525
$ws.Range("H113").Value = 2727
$ws.Range("I113").Value = 2727
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2727
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -557
$ws.Range("H122").Value = 5702.3706
$ws.Range("I122").Value = 5267.125
$ws.Range("J122").Value = 6335.4546
$ws.Range("K122").Value = 15801.375
$ws.Range("L122").Value = 19006.3638
$ws.Range("M122").Value = -13351.375
$ws.Range("N122").Value = -23906.3638
$ws.Range("H136").Value = 3151.78
$ws.Range("J136").Value = 3508.5334
$ws.Range("L136").Value = 10525.6002
$ws.Range("N136").Value = -15625.6002

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 11249.667
$ws.Range("J30").Value = 11249.667
$ws.Range("L30").Value = 11249.667
$ws.Range("N30").Value = -11463.667
$ws.Range("H39").Value = 26000
$ws.Range("J39").Value = 26000
$ws.Range("L39").Value = 26000
$ws.Range("N39").Value = -26826
$ws.Range("H122").Value = 64444.617
$ws.Range("I122").Value = 70738.94500000001
$ws.Range("K122").Value = 212216.835
$ws.Range("M122").Value = -209766.835
$ws.Range("H126").Value = 4255.9287
$ws.Range("I126").Value = 3960.5
$ws.Range("J126").Value = 4994.5
$ws.Range("K126").Value = 11881.5
$ws.Range("L126").Value = 14983.5
$ws.Range("M126").Value = -9411.5
$ws.Range("N126").Value = -19923.5
